# Insert a new weekly price record at row 480 of the "Zapallo" sheet,
# pushing the existing rows 480:525 down to 481:526.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(480).Insert()

$ws.Range("A480").Value = 4
$ws.Range("B480").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C480").Value = "Los Lagos"
$ws.Range("D480").Value = 45166
$ws.Range("E480").Value = 10
$ws.Range("F480").Value = 100112045
$ws.Range("G480").Value = "Zapallo"
$ws.Range("H480").Value = "Paine"
$ws.Range("I480").Value = "1a (guarda)"
$ws.Range("J480").Value = 250
$ws.Range("K480").Value = 650
$ws.Range("L480").Value = 650
$ws.Range("M480").Value = 650
$ws.Range("N480").Value = "$/kilo (volumen en unidades)"
$ws.Range("O480").Value = "Región de O'Higgins"
$ws.Range("P480").Value = 650
$ws.Range("Q480").Value = 1
$ws.Range("R480").Value = "Hortaliza"
